# Apply the commit's changes to the "swallowing_sillabe_rfe" results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: Validation -> F1 train
$ws.Range("O1").Value = "F1 train"

# O column (Validation/F1 train score) updates for rows 2-5
$ws.Range("O2").Value = 0.9295774647887324
$ws.Range("O3").Value = 0.9253731343283582
$ws.Range("O4").Value = 0.9722222222222222
$ws.Range("O5").Value = 1

# Row 6 (MLP, '5' technique) - parameters and metrics updated
$ws.Range("C6").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.7
$ws.Range("J6").Value = 0.7272727272727273
$ws.Range("K6").Value = 0.8888888888888888
$ws.Range("L6").Value = 0.6153846153846154
$ws.Range("M6").Value = 0.5454545454545454
$ws.Range("N6").Value = 0.8888888888888888
$ws.Range("O6").Value = 0.7088607594936709

# O column updates for rows 7-15
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 0.9411764705882353
$ws.Range("O9").Value = 0.9166666666666666
$ws.Range("O10").Value = 0.7575757575757576
$ws.Range("O11").Value = 0.6078431372549019
$ws.Range("O12").Value = 0.8055555555555556
$ws.Range("O13").Value = 1
$ws.Range("O14").Value = 1
$ws.Range("O15").Value = 0.8571428571428571

# Row 16 (MLP, 'Free' technique) - parameters and metrics updated
$ws.Range("C16").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 6
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 0.5
$ws.Range("J16").Value = 0.5454545454545454
$ws.Range("K16").Value = 0.6666666666666666
$ws.Range("L16").Value = 0.4615384615384616
$ws.Range("M16").Value = 0.3636363636363636
$ws.Range("N16").Value = 0.8888888888888888
$ws.Range("O16").Value = 0.7605633802816901
